$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "案由" (case reason) cell: "卖买合同" + "纠纷" (two runs) -> single run
#    "申请宣告公民限制民事行为能力"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("卖买合同纠纷", $true, $false, $false, $false, $false,
                         $true, 1, $false, "申请宣告公民限制民事行为能力", 2)

# ---------------------------------------------------------------------
# 2. "涉及标的额" (amount in dispute) cell: "266万" -> "0"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("266万", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0", 2)

# ---------------------------------------------------------------------
# 3. "实收费" (fee actually received) cell: drop the two leading
#    space-only runs, then change the amount "10000" -> "3000"
# ---------------------------------------------------------------------
$feeCell = $d.Tables(1).Rows(2).Cells(4)
$feeRange = $feeCell.Range
$leadingSpaces = $d.Range($feeRange.Start, $feeRange.Start + 2)
$leadingSpaces.Text = ""
$feeCell2 = $d.Tables(1).Rows(2).Cells(4)
$feeScoped = $d.Range($feeCell2.Range.Start, $feeCell2.Range.End)
$feeScoped.Find.Execute("10000", $true, $false, $false, $false, $false,
                         $true, 0, $false, "3000", 2)

# ---------------------------------------------------------------------
# 4. "委托人姓名" (client name) cell: "富卓汽车内饰（安徽）有限公司" -> "李桂芬"
#    Word also drops a "_GoBack" bookmark at the edit point (last-edit
#    position marker). We insert it right after the new text using a
#    throwaway trailing placeholder character so the zero-width bookmark
#    range lands correctly (then the placeholder is removed).
# ---------------------------------------------------------------------
$clientCell = $d.Tables(1).Rows(3).Cells(2)
$clientScoped = $d.Range($clientCell.Range.Start, $clientCell.Range.End)
$clientScoped.Find.Execute("富卓汽车内饰（安徽）有限公司", $true, $false, $false, $false, $false,
                            $true, 0, $false, "李桂芬#", 2)

$clientCell2 = $d.Tables(1).Rows(3).Cells(2)
$clientRange2 = $clientCell2.Range
$bmPos = $d.Range($clientRange2.End - 2, $clientRange2.End - 2)
$bmPos.Bookmarks.Add("_GoBack")

$clientCell3 = $d.Tables(1).Rows(3).Cells(2)
$clientScoped3 = $d.Range($clientCell3.Range.Start, $clientCell3.Range.End)
$clientScoped3.Find.Execute("#", $true, $false, $false, $false, $false,
                             $true, 0, $false, "", 2)

# ---------------------------------------------------------------------
# 5. "对方" (opposing party) cell: "芜湖亚利华汽车部件有限公司" -> "李忠楷"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("芜湖亚利华汽车部件有限公司", $true, $false, $false, $false, $false,
                         $true, 1, $false, "李忠楷", 2)

# ---------------------------------------------------------------------
# 6. "办案单位" (handling unit) cell: "芜湖市" -> "芜湖", "中级" -> "三山经济开发区人民"
#    (the third run "法院" is left untouched)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("芜湖市", $true, $false, $false, $false, $false,
                         $true, 1, $false, "芜湖", 2)
$d.Content.Find.Execute("中级", $true, $false, $false, $false, $false,
                         $true, 1, $false, "三山经济开发区人民", 2)
